$d = $word.ActiveDocument

# --- Hunk 1: merge the two timestamp runs ("Sat Nov 10" + " 11:30:33 PST 2017") into one run ---
$null = $d.Content.Find.Execute("Sat Nov 10 11:30:33 PST 2017", $false, $false, $false, $false, $false, $true, 1, $false, "Sat Nov 10 11:30:33 PST 2017", 2)

# --- Hunk 2: append the new "TUE Dec 19" purchase-details block after the
#     "Amount balance ... - 16361.0" paragraph, before the trailing blank paragraphs ---
$anchorIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "16361\.0") {
        $anchorIdx = $i
        break
    }
}
if ($anchorIdx -eq 0) { throw "anchor paragraph not found" }

$cur = $d.Paragraphs.Item($anchorIdx).Range
$cur.Collapse(0)

# --- new paragraph 1/15: bold=True color=None text='' ---
$cur.Text = "`r"

# --- new paragraph 2/15: bold=False color=None text='TUE Dec 19 10:55:05 PST 2017' ---
$cur.Text = "`rTUE Dec 19 10:55:05 PST 2017"
$newp = $d.Paragraphs.Item($anchorIdx + 2)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$newr.Font.Bold = 0
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 3/15: bold=False color=None text='Person Name\t\t\t\t- HSJ' ---
$cur.Text = "`rPerson Name`t`t`t`t- HSJ"
$newp = $d.Paragraphs.Item($anchorIdx + 3)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 4/15: bold=False color=None text='Bill number\t\t\t\t- 2054' ---
$cur.Text = "`rBill number`t`t`t`t- 2054"
$newp = $d.Paragraphs.Item($anchorIdx + 4)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 5/15: bold=False color=None text='---------------------------------------------------------------' ---
$cur.Text = "`r---------------------------------------------------------------"
$newp = $d.Paragraphs.Item($anchorIdx + 5)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 6/15: bold=False color=None text='Item Name\t\t\t\t- BEET' ---
$cur.Text = "`rItem Name`t`t`t`t- BEET"
$newp = $d.Paragraphs.Item($anchorIdx + 6)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 7/15: bold=False color=None text='Number of Pockets\t\t\t- 1' ---
$cur.Text = "`rNumber of Pockets`t`t`t- 1"
$newp = $d.Paragraphs.Item($anchorIdx + 7)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 8/15: bold=False color=None text='Number of KGs\t\t\t- 59' ---
$cur.Text = "`rNumber of KGs`t`t`t- 59"
$newp = $d.Paragraphs.Item($anchorIdx + 8)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 9/15: bold=False color=None text='Rate\t\t\t\t\t- 22' ---
$cur.Text = "`rRate`t`t`t`t`t- 22"
$newp = $d.Paragraphs.Item($anchorIdx + 9)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 10/15: bold=False color=None text='Total Price\t\t\t\t- 1298.0' ---
$cur.Text = "`rTotal Price`t`t`t`t- 1298.0"
$newp = $d.Paragraphs.Item($anchorIdx + 10)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 11/15: bold=False color=FF0000 text='Amount Received\t\t\t- 6000' ---
$cur.Text = "`rAmount Received`t`t`t- 6000"
$newp = $d.Paragraphs.Item($anchorIdx + 11)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$newr.Font.Color = 255
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 12/15: bold=True color=None text='Amount balance\t\t\t- 11659.0' ---
$cur.Text = "`rAmount balance`t`t`t- 11659.0"
$newp = $d.Paragraphs.Item($anchorIdx + 12)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$newr.Font.Bold = 1
$newr.Font.Color = -16777216
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 13/15: bold=False color=None text='Amount Received mode\t\t- CASH' ---
$cur.Text = "`rAmount Received mode`t`t- CASH"
$newp = $d.Paragraphs.Item($anchorIdx + 13)
$newr = $newp.Range
$newr.Font.NameAscii = "Courier New"
$newr.Font.NameOther = "Courier New"
$newr.Font.NameBi = "Courier New"
$newr.Font.Bold = 0
$cur = $newp.Range
$cur.Collapse(0)

# --- new paragraph 14/15: bold=False color=None text='' ---
$cur.Text = "`r"

# --- new paragraph 15/15: bold=True color=None text='' ---
$cur.Font.Bold = 1
$cur.Text = "`r"

Write-Output "anchorIdx=$anchorIdx paragraphCount=$($d.Paragraphs.Count)"